# Update "ValueSet-snomed-myelodysplastic-diseases-vs.xlsx" metadata:
# new version, status, date, contact info, and a new "Jurisdiction" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Simple value updates on the Metadata sheet ---
$ws.Range("B3").Value2  = "0.1.7"
$ws.Range("B6").Value2  = "draft"
$ws.Range("B8").Value2  = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after the two Contact rows (row 12) ---
$ws.Rows.Item(12).Insert()

# Match the formatting of the surrounding data rows.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A12").Value2 = "Jurisdiction"
$ws.Range("B12").Value2 = ""
